# Apply the latest cryptos snapshot values scraped for this update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cell reference, new text value) pairs, in sheet order.
$updates = @(
    @("D2", "27.454.80"),
    @("E2", "  -0.95%  "),
    @("D3", "1.829.32"),
    @("E3", "  -2.06%  "),
    @("E4", "  -0.91%  "),
    @("E5", "  -1.01%  "),
    @("E6", "  -0.76%  "),
    @("D7", "0.4573"),
    @("E7", "  -2.44%  "),
    @("E8", "  -2.83%  "),
    @("D9", "46.31"),
    @("E9", "  +1.42%  "),
    @("D10", "0.07879"),
    @("E10", "  -1.51%  "),
    @("D11", "0.9661"),
    @("E11", "  -3.92%  "),
    @("D12", "21.00"),
    @("E12", "  -3.91%  "),
    @("B13", "WrappedEther"),
    @("C13", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"),
    @("D13", "1.829.64"),
    @("E13", "  -2.29%  "),
    @("B14", "Polkadot"),
    @("C14", "https://coinranking.com/coin/25W7FG7om+polkadot-dot"),
    @("D14", "5.870"),
    @("E14", "  -2.13%  "),
    @("D15", "7.082"),
    @("E15", "  -2.37%  "),
    @("D16", "1.003"),
    @("E16", "  -0.82%  "),
    @("E17", "  +0.80%  "),
    @("D18", "0.06600"),
    @("E18", "  -2.23%  "),
    @("D19", "0.00001025"),
    @("E19", "  -1.79%  "),
    @("D20", "17.08"),
    @("E20", "  -0.73%  "),
    @("E21", "  -0.75%  "),
    @("D22", "27.449.79"),
    @("E22", "  -0.99%  "),
    @("D23", "5.329"),
    @("E23", "  -2.69%  "),
    @("E24", "  -1.27%  "),
    @("D25", "2.287"),
    @("E25", "  -1.36%  "),
    @("D26", "2.048.77"),
    @("E26", "  -2.14%  "),
    @("D27", "155.65"),
    @("E27", "  -1.90%  "),
    @("D28", "19.35"),
    @("E28", "  -2.26%  "),
    @("D29", "2.068"),
    @("E29", "  -4.19%  "),
    @("D30", "5.288"),
    @("E30", "  -2.96%  "),
    @("D31", "118.41"),
    @("E31", "  -2.81%  "),
    @("D32", "0.9391"),
    @("E32", "  -4.16%  "),
    @("D33", "0.09291"),
    @("E33", "  -2.05%  "),
    @("D34", "3.577"),
    @("E34", "  -1.31%  "),
    @("D35", "5.234"),
    @("E35", "  -1.82%  "),
    @("D36", "1.325"),
    @("E36", "  -1.25%  "),
    @("D37", "0.05919"),
    @("E37", "  -2.27%  "),
    @("D38", "0.02184"),
    @("E38", "  -2.33%  "),
    @("D39", "8.105"),
    @("E39", "  -2.23%  "),
    @("D40", "1.148"),
    @("E40", "  -4.02%  "),
    @("D41", "0.5768"),
    @("E41", "  -3.50%  "),
    @("D42", "0.1825"),
    @("E42", "  -3.29%  "),
    @("D43", "9.977"),
    @("E43", "  -3.30%  "),
    @("E44", "  +1.88%  "),
    @("D45", "11.96"),
    @("E45", "  -1.95%  "),
    @("D46", "0.5439"),
    @("E46", "  -3.92%  "),
    @("D47", "1.867"),
    @("E47", "  -3.03%  "),
    @("D48", "110.47"),
    @("E48", "  -1.40%  "),
    @("D49", "0.06586"),
    @("E49", "  -2.59%  "),
    @("D50", "1.001"),
    @("D51", "1.042"),
    @("E51", "  -1.41%  ")
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $range = $ws.Range($cellRef)
    if ($cellRef.StartsWith("D")) {
        # Price column values are stored as literal text (e.g. "27.454.80",
        # "5.870"); force text format first so Excel does not silently
        # reinterpret numeric-looking text as a Number (which would also
        # drop significant trailing zeros, e.g. "5.870" -> 5.87).
        $range.NumberFormat = "@"
    }
    $range.Value = $newValue
}
